# chore: update Sheets via scheduled runner
# Refreshes market-price / profit figures (cols H-N) on a handful of
# leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 306.2
$ws.Range("I42").Value = 84
$ws.Range("J42").Value = 361.75
$ws.Range("K42").Value = 252
$ws.Range("L42").Value = 1085.25
$ws.Range("M42").Value = -22
$ws.Range("N42").Value = -1545.25

$ws.Range("H112").Value = 2882.1428
$ws.Range("J112").Value = 2918.5186
$ws.Range("L112").Value = 8755.5558
$ws.Range("N112").Value = -10971.5558

$ws.Range("H129").Value = 1230.9667
$ws.Range("J129").Value = 1505.65
$ws.Range("L129").Value = 4516.950000000001
$ws.Range("N129").Value = -14516.95

$ws.Range("H132").Value = 5846.1816
$ws.Range("I132").Value = 6038.3687
$ws.Range("K132").Value = 18115.1061
$ws.Range("M132").Value = -15585.1061

$ws.Range("H138").Value = 2752.8628
$ws.Range("I138").Value = 3809
$ws.Range("J138").Value = 2495.2683
$ws.Range("K138").Value = 11427
$ws.Range("L138").Value = 7485.804900000001
$ws.Range("M138").Value = -6287
$ws.Range("N138").Value = -17765.8049

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11907338
$ws.Range("I61").Value = 41668010
$ws.Range("J61").Value = 3069.6
$ws.Range("K61").Value = 41668010
$ws.Range("L61").Value = 3069.6
$ws.Range("M61").Value = -41667798
$ws.Range("N61").Value = -3493.6

$ws.Range("H74").Value = 1670.5625
$ws.Range("I74").Value = 1266.8334
$ws.Range("J74").Value = 1912.8
$ws.Range("K74").Value = 1266.8334
$ws.Range("L74").Value = 1912.8
$ws.Range("M74").Value = -392.8334
$ws.Range("N74").Value = -3660.8

$ws.Range("H77").Value = 1670.5625
$ws.Range("I77").Value = 1266.8334
$ws.Range("J77").Value = 1912.8
$ws.Range("K77").Value = 6334.166999999999
$ws.Range("L77").Value = 9564
$ws.Range("M77").Value = -1966.166999999999
$ws.Range("N77").Value = -18300

$ws.Range("H123").Value = 29814.916
$ws.Range("J123").Value = 29814.916
$ws.Range("L123").Value = 29814.916
$ws.Range("N123").Value = -39614.916

$ws.Range("H136").Value = 11907338
$ws.Range("I136").Value = 41668010
$ws.Range("J136").Value = 3069.6
$ws.Range("K136").Value = 125004030
$ws.Range("L136").Value = 9208.799999999999
$ws.Range("M136").Value = -125001480
$ws.Range("N136").Value = -14308.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H131").Value = 99780
$ws.Range("J131").Value = 99780
$ws.Range("L131").Value = 99780
$ws.Range("N131").Value = -109860

$ws.Range("H134").Value = 2863.543
$ws.Range("I134").Value = 2686.7827
$ws.Range("K134").Value = 8060.348100000001
$ws.Range("M134").Value = -5525.348100000001

$ws.Range("H141").Value = 19625
$ws.Range("J141").Value = 19625
$ws.Range("L141").Value = 19625
$ws.Range("N141").Value = -29985

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4911.2705
$ws.Range("I31").Value = 1450.0555
$ws.Range("K31").Value = 1450.0555
$ws.Range("M31").Value = -1155.0555

$ws.Range("H34").Value = 4911.2705
$ws.Range("I34").Value = 1450.0555
$ws.Range("K34").Value = 1450.0555
$ws.Range("M34").Value = -1248.0555

$ws.Range("H58").Value = 2287.1333
$ws.Range("J58").Value = 2562.8
$ws.Range("L58").Value = 2562.8
$ws.Range("N58").Value = -2968.8

$ws.Range("H136").Value = 2287.1333
$ws.Range("J136").Value = 2562.8
$ws.Range("L136").Value = 7688.400000000001
$ws.Range("N136").Value = -12788.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1305.5186
$ws.Range("I129").Value = 642.1111
$ws.Range("J129").Value = 1637.2222
$ws.Range("K129").Value = 1926.3333
$ws.Range("L129").Value = 4911.6666
$ws.Range("M129").Value = 3073.6667
$ws.Range("N129").Value = -14911.6666

$ws.Range("H138").Value = 1701.2858
$ws.Range("J138").Value = 3141
$ws.Range("L138").Value = 9423
$ws.Range("N138").Value = -19703

$ws.Range("H139").Value = 5301.6875
$ws.Range("J139").Value = 21666.334
$ws.Range("L139").Value = 64999.00199999999
$ws.Range("N139").Value = -75279.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 80110
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 80110
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 80110
$ws.Range("N59").Value = -81276
$ws.Range("M59").ClearContents()

$ws.Range("H93").Value = 45597.867
$ws.Range("J93").Value = 45597.867
$ws.Range("L93").Value = 45597.867
$ws.Range("N93").Value = -49341.867

$ws.Range("H112").Value = 98293
$ws.Range("J112").Value = 98293
$ws.Range("L112").Value = 98293
$ws.Range("N112").Value = -100509

$ws.Range("H113").Value = 2088.875
$ws.Range("I113").Value = 1802.75
$ws.Range("J113").Value = 2375
$ws.Range("K113").Value = 1802.75
$ws.Range("L113").Value = 2375
$ws.Range("M113").Value = 367.25
$ws.Range("N113").Value = -6715

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H122").Value = 4986.6875
$ws.Range("I122").Value = 4467.8335
$ws.Range("J122").Value = 5298
$ws.Range("K122").Value = 13403.5005
$ws.Range("L122").Value = 15894
$ws.Range("M122").Value = -10953.5005
$ws.Range("N122").Value = -20794

$ws.Range("H132").Value = 2805.0588
$ws.Range("I132").Value = 2790.8572
$ws.Range("J132").Value = 2815
$ws.Range("K132").Value = 8372.571599999999
$ws.Range("L132").Value = 8445
$ws.Range("M132").Value = -5842.571599999999
$ws.Range("N132").Value = -13505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7905.1333
$ws.Range("I22").Value = 760
$ws.Range("J22").Value = 22195.4
$ws.Range("K22").Value = 760
$ws.Range("L22").Value = 22195.4
$ws.Range("M22").Value = -465
$ws.Range("N22").Value = -22785.4

$ws.Range("H27").Value = 7905.1333
$ws.Range("I27").Value = 760
$ws.Range("J27").Value = 22195.4
$ws.Range("K27").Value = 760
$ws.Range("L27").Value = 22195.4
$ws.Range("M27").Value = -653
$ws.Range("N27").Value = -22409.4

$ws.Range("H132").Value = 2919.9756
$ws.Range("I132").Value = 2512.4075
$ws.Range("J132").Value = 3706
$ws.Range("K132").Value = 7537.2225
$ws.Range("L132").Value = 11118
$ws.Range("M132").Value = -5007.2225
$ws.Range("N132").Value = -16178

$ws.Range("H136").Value = 3268867.2
$ws.Range("I136").Value = 1004.1071
$ws.Range("J136").Value = 7247135.5
$ws.Range("K136").Value = 3012.3213
$ws.Range("L136").Value = 21741406.5
$ws.Range("M136").Value = -462.3212999999996
$ws.Range("N136").Value = -21746506.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 25031.223
$ws.Range("J123").Value = 25031.223
$ws.Range("L123").Value = 25031.223
$ws.Range("N123").Value = -34831.223

$ws.Range("H136").Value = 2404.8135
$ws.Range("I136").Value = 2124.2222
$ws.Range("J136").Value = 3306.7144
$ws.Range("K136").Value = 6372.6666
$ws.Range("L136").Value = 9920.143199999999
$ws.Range("M136").Value = -3822.6666
$ws.Range("N136").Value = -15020.1432
